$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_WVR = $wb.Worksheets.Item("WVR")

# Row 88 (ALC) - diff hunk @@ -5046,25 +5046,25 @@
$ws_ALC.Range("H88").Value = 1849.4286
$ws_ALC.Range("I88").Value = 2318.2
$ws_ALC.Range("J88").Value = 1589
$ws_ALC.Range("K88").Value = 2318.2
$ws_ALC.Range("L88").Value = 1589
$ws_ALC.Range("M88").Value = -1912.2
$ws_ALC.Range("N88").Value = -2401

# Row 91 (ALC) - diff hunk @@ -5202,25 +5202,25 @@
$ws_ALC.Range("H91").Value = 1849.4286
$ws_ALC.Range("I91").Value = 2318.2
$ws_ALC.Range("J91").Value = 1589
$ws_ALC.Range("K91").Value = 2318.2
$ws_ALC.Range("L91").Value = 1589
$ws_ALC.Range("M91").Value = -914.1999999999998
$ws_ALC.Range("N91").Value = -4397

# Row 127 (ALC) - diff hunk @@ -6999,25 +6999,25 @@
$ws_ALC.Range("H127").Value = 2899.8
$ws_ALC.Range("I127").Value = 3250
$ws_ALC.Range("J127").Value = 1499
$ws_ALC.Range("K127").Value = 9750
$ws_ALC.Range("L127").Value = 4497
$ws_ALC.Range("M127").Value = -4790
$ws_ALC.Range("N127").Value = -14417

# Row 137 (ALC) - diff hunk @@ -7495,22 +7495,22 @@
$ws_ALC.Range("H137").Value = 1277070.8
$ws_ALC.Range("I137").Value = 993299.7
$ws_ALC.Range("J137").Value = 1442603.8
$ws_ALC.Range("K137").Value = 2979899.1
$ws_ALC.Range("L137").Value = 4327811.4
$ws_ALC.Range("M137").Value = -2977349.1

# Row 32 (ARM) - diff hunk @@ -9334,25 +9334,22 @@
$ws_ARM.Range("H32").Value = 13999.333
$ws_ARM.Range("I32").Value = 0
$ws_ARM.Range("J32").Value = 13999.333
$ws_ARM.Range("K32").Value = 0
$ws_ARM.Range("L32").Value = 13999.333
$ws_ARM.Range("N32").Value = -14573.333
$ws_ARM.Range("M32").ClearContents()

# Row 45 (ARM) - diff hunk @@ -9965,22 +9962,22 @@
$ws_ARM.Range("H45").Value = 4586.1113
$ws_ARM.Range("I45").Value = 3797.0588
$ws_ARM.Range("J45").Value = 18000
$ws_ARM.Range("K45").Value = 3797.0588
$ws_ARM.Range("L45").Value = 18000
$ws_ARM.Range("M45").Value = -3420.0588

# Row 74 (ARM) - diff hunk @@ -11395,25 +11392,25 @@
$ws_ARM.Range("H74").Value = 2317793.8
$ws_ARM.Range("I74").Value = 2843685
$ws_ARM.Range("J74").Value = 3872.3
$ws_ARM.Range("K74").Value = 2843685
$ws_ARM.Range("L74").Value = 3872.3
$ws_ARM.Range("M74").Value = -2842811
$ws_ARM.Range("N74").Value = -5620.3

# Row 77 (ARM) - diff hunk @@ -11545,25 +11542,25 @@
$ws_ARM.Range("H77").Value = 2317793.8
$ws_ARM.Range("I77").Value = 2843685
$ws_ARM.Range("J77").Value = 3872.3
$ws_ARM.Range("K77").Value = 14218425
$ws_ARM.Range("L77").Value = 19361.5
$ws_ARM.Range("M77").Value = -14214057
$ws_ARM.Range("N77").Value = -28097.5

# Row 97 (ARM) - diff hunk @@ -12522,22 +12519,22 @@
$ws_ARM.Range("H97").Value = 1316.7142
$ws_ARM.Range("I97").Value = 1369.5
$ws_ARM.Range("J97").Value = 1000
$ws_ARM.Range("K97").Value = 1369.5
$ws_ARM.Range("L97").Value = 1000
$ws_ARM.Range("M97").Value = -873.5

# Row 2 (BSM) - diff hunk @@ -14812,22 +14809,22 @@
$ws_BSM.Range("H2").Value = 79916
$ws_BSM.Range("I2").Value = 0
$ws_BSM.Range("J2").Value = 79916
$ws_BSM.Range("K2").Value = 0
$ws_BSM.Range("L2").Value = 79916
$ws_BSM.Range("N2").Value = -80142

# Row 20 (BSM) - diff hunk @@ -15706,25 +15703,25 @@
$ws_BSM.Range("H20").Value = 2125.2666
$ws_BSM.Range("I20").Value = 1857.4762
$ws_BSM.Range("J20").Value = 2750.111
$ws_BSM.Range("K20").Value = 1857.4762
$ws_BSM.Range("L20").Value = 2750.111
$ws_BSM.Range("M20").Value = -1610.4762
$ws_BSM.Range("N20").Value = -3244.111

# Row 94 (BSM) - diff hunk @@ -19320,22 +19317,22 @@
$ws_BSM.Range("H94").Value = 457.42105
$ws_BSM.Range("I94").Value = 448.1111
$ws_BSM.Range("J94").Value = 625
$ws_BSM.Range("K94").Value = 448.1111
$ws_BSM.Range("L94").Value = 625
$ws_BSM.Range("M94").Value = 2.888899999999978

# Row 105 (BSM) - diff hunk @@ -19871,25 +19868,25 @@
$ws_BSM.Range("H105").Value = 3491.64
$ws_BSM.Range("I105").Value = 2510.125
$ws_BSM.Range("J105").Value = 5236.5557
$ws_BSM.Range("K105").Value = 2510.125
$ws_BSM.Range("L105").Value = 5236.5557
$ws_BSM.Range("M105").Value = -763.125
$ws_BSM.Range("N105").Value = -8730.555700000001

# Row 134 (BSM) - diff hunk @@ -21283,22 +21280,22 @@
$ws_BSM.Range("H134").Value = 599611.25
$ws_BSM.Range("I134").Value = 885789.0600000001
$ws_BSM.Range("J134").Value = 5241.923
$ws_BSM.Range("K134").Value = 2657367.18
$ws_BSM.Range("L134").Value = 15725.769
$ws_BSM.Range("M134").Value = -2654832.18

# Row 137 (BSM) - diff hunk @@ -21430,25 +21427,25 @@
$ws_BSM.Range("H137").Value = 98836
$ws_BSM.Range("I137").Value = 77000
$ws_BSM.Range("J137").Value = 113393.336
$ws_BSM.Range("K137").Value = 77000
$ws_BSM.Range("L137").Value = 113393.336
$ws_BSM.Range("M137").Value = -71900
$ws_BSM.Range("N137").Value = -123593.336

# Row 31 (CRP) - diff hunk @@ -23205,22 +23202,22 @@
$ws_CRP.Range("H31").Value = 8631.027
$ws_CRP.Range("I31").Value = 2755.5386
$ws_CRP.Range("J31").Value = 11951.956
$ws_CRP.Range("K31").Value = 2755.5386
$ws_CRP.Range("L31").Value = 11951.956
$ws_CRP.Range("M31").Value = -2460.5386

# Row 34 (CRP) - diff hunk @@ -23352,22 +23349,22 @@
$ws_CRP.Range("H34").Value = 8631.027
$ws_CRP.Range("I34").Value = 2755.5386
$ws_CRP.Range("J34").Value = 11951.956
$ws_CRP.Range("K34").Value = 2755.5386
$ws_CRP.Range("L34").Value = 11951.956
$ws_CRP.Range("M34").Value = -2553.5386

# Row 70 (CRP) - diff hunk @@ -25110,19 +25107,22 @@
$ws_CRP.Range("H70").Value = 67602.75
$ws_CRP.Range("I70").Value = 0
$ws_CRP.Range("J70").Value = 67602.75
$ws_CRP.Range("K70").Value = 0
$ws_CRP.Range("L70").Value = 67602.75
$ws_CRP.Range("N70").Value = -68232.75

# Row 73 (CRP) - diff hunk @@ -25257,19 +25257,22 @@
$ws_CRP.Range("H73").Value = 67602.75
$ws_CRP.Range("I73").Value = 0
$ws_CRP.Range("J73").Value = 67602.75
$ws_CRP.Range("K73").Value = 0
$ws_CRP.Range("L73").Value = 67602.75
$ws_CRP.Range("N73").Value = -69786.75

# Row 80 (CRP) - diff hunk @@ -25603,19 +25606,22 @@
$ws_CRP.Range("H80").Value = 52776
$ws_CRP.Range("I80").Value = 0
$ws_CRP.Range("J80").Value = 52776
$ws_CRP.Range("K80").Value = 0
$ws_CRP.Range("L80").Value = 52776
$ws_CRP.Range("N80").Value = -55022

# Row 83 (CRP) - diff hunk @@ -25744,19 +25750,22 @@
$ws_CRP.Range("H83").Value = 52776
$ws_CRP.Range("I83").Value = 0
$ws_CRP.Range("J83").Value = 52776
$ws_CRP.Range("K83").Value = 0
$ws_CRP.Range("L83").Value = 158328
$ws_CRP.Range("N83").Value = -169560

# Row 97 (CRP) - diff hunk @@ -26433,22 +26442,22 @@
$ws_CRP.Range("H97").Value = 60000
$ws_CRP.Range("I97").Value = 0
$ws_CRP.Range("J97").Value = 60000
$ws_CRP.Range("K97").Value = 0
$ws_CRP.Range("L97").Value = 60000
$ws_CRP.Range("N97").Value = -61982

# Row 134 (CRP) - diff hunk @@ -28222,22 +28231,22 @@
$ws_CRP.Range("H134").Value = 3177.4375
$ws_CRP.Range("I134").Value = 2161.4827
$ws_CRP.Range("J134").Value = 12998.333
$ws_CRP.Range("K134").Value = 6484.4481
$ws_CRP.Range("L134").Value = 38994.999
$ws_CRP.Range("M134").Value = -3949.4481

# Row 25 (CUL) - diff hunk @@ -29910,22 +29919,19 @@
$ws_CUL.Range("H25").Value = 0
$ws_CUL.Range("I25").Value = 0
$ws_CUL.Range("J25").Value = 0
$ws_CUL.Range("K25").Value = 0
$ws_CUL.Range("L25").Value = 0
$ws_CUL.Range("M25").ClearContents()

# Row 30 (CUL) - diff hunk @@ -30164,22 +30170,19 @@
$ws_CUL.Range("H30").Value = 0
$ws_CUL.Range("I30").Value = 0
$ws_CUL.Range("J30").Value = 0
$ws_CUL.Range("K30").Value = 0
$ws_CUL.Range("L30").Value = 0
$ws_CUL.Range("M30").ClearContents()

# Row 131 (CUL) - diff hunk @@ -35302,25 +35305,25 @@
$ws_CUL.Range("H131").Value = 11003.692
$ws_CUL.Range("I131").Value = 4299
$ws_CUL.Range("J131").Value = 13015.1
$ws_CUL.Range("K131").Value = 12897
$ws_CUL.Range("L131").Value = 39045.3
$ws_CUL.Range("M131").Value = -7857
$ws_CUL.Range("N131").Value = -49125.3

# Row 15 (GSM) - diff hunk @@ -36602,25 +36605,22 @@
$ws_GSM.Range("H15").Value = 41909.09
$ws_GSM.Range("I15").Value = 0
$ws_GSM.Range("J15").Value = 41909.09
$ws_GSM.Range("K15").Value = 0
$ws_GSM.Range("L15").Value = 41909.09
$ws_GSM.Range("N15").Value = -42485.09
$ws_GSM.Range("M15").ClearContents()

# Row 81 (GSM) - diff hunk @@ -39824,25 +39824,22 @@
$ws_GSM.Range("H81").Value = 41909.09
$ws_GSM.Range("I81").Value = 0
$ws_GSM.Range("J81").Value = 41909.09
$ws_GSM.Range("K81").Value = 0
$ws_GSM.Range("L81").Value = 41909.09
$ws_GSM.Range("N81").Value = -43905.09
$ws_GSM.Range("M81").ClearContents()

# Row 84 (GSM) - diff hunk @@ -39980,25 +39977,22 @@
$ws_GSM.Range("H84").Value = 41909.09
$ws_GSM.Range("I84").Value = 0
$ws_GSM.Range("J84").Value = 41909.09
$ws_GSM.Range("K84").Value = 0
$ws_GSM.Range("L84").Value = 125727.27
$ws_GSM.Range("N84").Value = -135711.27
$ws_GSM.Range("M84").ClearContents()

# Row 97 (GSM) - diff hunk @@ -40611,25 +40605,25 @@
$ws_GSM.Range("H97").Value = 1077.963
$ws_GSM.Range("I97").Value = 362.66666
$ws_GSM.Range("J97").Value = 1972.0834
$ws_GSM.Range("K97").Value = 362.66666
$ws_GSM.Range("L97").Value = 1972.0834
$ws_GSM.Range("M97").Value = 133.33334
$ws_GSM.Range("N97").Value = -2964.0834

# Row 98 (GSM) - diff hunk @@ -40663,22 +40657,22 @@
$ws_GSM.Range("H98").Value = 73367.5
$ws_GSM.Range("I98").Value = 0
$ws_GSM.Range("J98").Value = 73367.5
$ws_GSM.Range("K98").Value = 0
$ws_GSM.Range("L98").Value = 73367.5
$ws_GSM.Range("N98").Value = -79357.5

# Row 102 (GSM) - diff hunk @@ -40862,25 +40856,25 @@
$ws_GSM.Range("H102").Value = 2402.9575
$ws_GSM.Range("I102").Value = 1684.3823
$ws_GSM.Range("J102").Value = 4282.3076
$ws_GSM.Range("K102").Value = 1684.3823
$ws_GSM.Range("L102").Value = 4282.3076
$ws_GSM.Range("M102").Value = -62.38229999999999
$ws_GSM.Range("N102").Value = -7526.3076

# Row 126 (GSM) - diff hunk @@ -42026,22 +42020,22 @@
$ws_GSM.Range("H126").Value = 698781.3
$ws_GSM.Range("I126").Value = 1517251.5
$ws_GSM.Range("J126").Value = 6229.615
$ws_GSM.Range("K126").Value = 4551754.5
$ws_GSM.Range("L126").Value = 18688.845
$ws_GSM.Range("M126").Value = -4549284.5

# Row 136 (GSM) - diff hunk @@ -42513,22 +42507,22 @@
$ws_GSM.Range("H136").Value = 20692.861
$ws_GSM.Range("I136").Value = 0
$ws_GSM.Range("J136").Value = 20692.861
$ws_GSM.Range("K136").Value = 0
$ws_GSM.Range("L136").Value = 62078.583
$ws_GSM.Range("N136").Value = -67178.583

# Row 132 (LTW) - diff hunk @@ -49301,25 +49295,25 @@
$ws_LTW.Range("H132").Value = 723833.75
$ws_LTW.Range("I132").Value = 937406
$ws_LTW.Range("J132").Value = 5454.5454
$ws_LTW.Range("K132").Value = 2812218
$ws_LTW.Range("L132").Value = 16363.6362
$ws_LTW.Range("M132").Value = -2809688
$ws_LTW.Range("N132").Value = -21423.6362

# Row 136 (LTW) - diff hunk @@ -49500,25 +49494,25 @@
$ws_LTW.Range("H136").Value = 4521.9585
$ws_LTW.Range("I136").Value = 3618.4358
$ws_LTW.Range("J136").Value = 8437.223
$ws_LTW.Range("K136").Value = 10855.3074
$ws_LTW.Range("L136").Value = 25311.669
$ws_LTW.Range("M136").Value = -8305.307400000002
$ws_LTW.Range("N136").Value = -30411.669

# Row 81 (WVR) - diff hunk @@ -53741,22 +53735,22 @@
$ws_WVR.Range("H81").Value = 1012.5333
$ws_WVR.Range("I81").Value = 969
$ws_WVR.Range("J81").Value = 1099.6
$ws_WVR.Range("K81").Value = 1938
$ws_WVR.Range("L81").Value = 2199.2
$ws_WVR.Range("M81").Value = -877

# Row 84 (WVR) - diff hunk @@ -53891,22 +53885,22 @@
$ws_WVR.Range("H84").Value = 1012.5333
$ws_WVR.Range("I84").Value = 969
$ws_WVR.Range("J84").Value = 1099.6
$ws_WVR.Range("K84").Value = 9690
$ws_WVR.Range("L84").Value = 10996
$ws_WVR.Range("M84").Value = -4386

# Row 122 (WVR) - diff hunk @@ -55735,22 +55729,22 @@
$ws_WVR.Range("H122").Value = 1969.5814
$ws_WVR.Range("I122").Value = 1491.0571
$ws_WVR.Range("J122").Value = 4063.125
$ws_WVR.Range("K122").Value = 4473.1713
$ws_WVR.Range("L122").Value = 12189.375
$ws_WVR.Range("M122").Value = -2023.1713
